$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A12" = -21.55939999999999
    "A27" = -21.68449999999999
    "A32" = -21.25
    "A36" = -20.0184
    "A38" = -19.3715
    "A46" = -21.4255
    "A54" = -21.59359999999999
    "A55" = -22.4374
    "A56" = -22.3601
    "A67" = -21.52169999999998
    "A69" = -21.63269999999998
    "A72" = -21.53169999999999
    "A83" = -21.7822
    "A86" = -22.26780000000002
    "A91" = -21.4683
    "A93" = -21.30449999999999
    "A99" = -20.26309999999999
    "C3"  = -12.2782
    "C14" = -13.5346
    "C26" = -12.26230000000001
    "C31" = -12.7181
    "C35" = -12.79920000000001
    "C37" = -13.5253
    "C45" = -13.95659999999999
    "C52" = -11.0246
    "C57" = -14.34399999999999
    "C81" = -12.8336
    "C83" = -12.78239999999999
    "C100" = -12.80219999999999
    "C102" = -13.1437
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
